# Update "Datos actualizados" timestamp footer
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 06:29"

# --- Row 7: India (no position change, values refreshed) ---
$ws.Range("B7").Value = 473105
$ws.Range("C7").Value = 120
$ws.Range("D7").Value = 271697
$ws.Range("E7").Value = 186501

# --- Rows 16/17: Pakistan overtakes Turquia in the ranking ---
# Row 16 becomes Pakistan with its freshly updated figures
$ws.Range("A16").Value = "Pakistan"
$ws.Range("B16").Value = 192970
$ws.Range("C16").Value = 4044
$ws.Range("D16").Value = 81307
$ws.Range("E16").Value = 107760
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 148
$ws.Range("H16").Value = 3903

# Row 17 becomes Turquia, carrying its previous (unchanged) figures
$ws.Range("A17").Value = "Turquia"
$ws.Range("B17").Value = 191657
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 164234
$ws.Range("E17").Value = 22398
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 5025

# --- Row 54: Kazajistan (no position change, values refreshed) ---
$ws.Range("B54").Value = 19285
$ws.Range("C54").Value = 520
$ws.Range("E54").Value = 7267

# --- Row 74: Australia (no position change, values refreshed) ---
$ws.Range("B74").Value = 7558
$ws.Range("C74").Value = 37
$ws.Range("D74").Value = 6931

# --- Row 92: Kirguistan (no position change, values refreshed) ---
$ws.Range("B92").Value = 3954
$ws.Range("C92").Value = 228
$ws.Range("D92").Value = 2112
$ws.Range("E92").Value = 1799
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 43

# --- Row 157: Surinam (no position change, values refreshed) ---
$ws.Range("E157").Value = 193
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 10

# --- Row 165: Mongolia (no position change, values refreshed) ---
$ws.Range("B165").Value = 216
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 169

# --- Rows 211/212: Montserrat overtakes Seychelles in the ranking ---
# Row 211 becomes Montserrat, carrying its previous (unchanged) figures
$ws.Range("A211").Value = "Montserrat"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1

# Row 212 becomes Seychelles, carrying its previous (unchanged) figures
$ws.Range("A212").Value = "Seychelles"
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 11
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

# --- Rows 208/209: Groenlandia overtakes Islas Malvinas (tie in data, only order swaps) ---
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"
